# Sun Jan 15 23:35:19 UTC 2023 - "Updated symbol list" GitHub Actions run.
# Refreshes the Price (column D) and Volume(1h) (column E) readings that
# moved since the previous scrape. Every cell in these columns is stored
# as literal text ("300.77", "-1.26%", ...), so each new value is written
# through Formula with a leading apostrophe - this is how COM/Excel force
# literal text instead of auto-coercing the string into a Number or
# Percentage cell. Style is reset to Normal right after so we don't leave
# a stray quote-prefix format behind (the source cells carry no style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'300.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'-1.26%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'31.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'-2.67%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Formula = "'5.137"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Formula = "'-3.10%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'0.07395"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'-2.52%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Formula = "'2.177"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'31.63%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Formula = "'7.923"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Formula = "'0.56%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Formula = "'3.823"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Formula = "'-0.98%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Formula = "'0.9185"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "'-1.12%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Formula = "'0.1709"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Formula = "'1.07%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Formula = "'0.07559"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "'-4.95%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Formula = "'0.08148"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "'0.83%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Formula = "'0.03013"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Formula = "'-1.29%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Formula = "'0.09928"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'-0.09%"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Formula = "'-0.55%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Formula = "'0.006074"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'-3.28%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Formula = "'0.84%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Formula = "'2.224"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "'-0.52%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Formula = "'-1.06%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Formula = "'-1.88%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Formula = "'4.648"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'2.27%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Formula = "'0.04637"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "'1.14%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Formula = "'0.1566"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'-3.04%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Formula = "'0.001225"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'0.99%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Formula = "'0.004475"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "'-0.29%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Formula = "'-6.98%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Formula = "'0.0003425"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Formula = "'101.86%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E39").Formula = "'0.61%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Formula = "'0.04514"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'-0.56%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Formula = "'0.007315"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Formula = "'0.1348"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Formula = "'-1.01%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Formula = "'0.002227"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'6.89%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Formula = "'-23.39%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Formula = "'0.00006270"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'2.08%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Formula = "'-22.87%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Formula = "'0.8085"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'14.02%"
$ws.Range("E47").Style = "Normal"
